$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the emailId value (A2) and the password value (B2)
$ws.Range("A2").Value = "Test_UX09@westpharma.com"
$ws.Range("B2").Value = "Westpharm@2019"

# B2 now also gets a hyperlink (mailto:) like A2 already had, with the
# matching "Hyperlink" cell style so it renders/underlines the same way.
$ws.Hyperlinks.Add($ws.Cells.Item(2, 2), "mailto:Westpharm@2019")
$ws.Range("B2").Style = "Hyperlink"

# Column B widens to fit the new text (closest attainable width to the
# recorded best-fit value of 16.54296875 given this engine's column-width
# quantization).
$ws.Columns("B").ColumnWidth = 15.6

# Selection moves to C2
$ws.Range("C2").Select()
